$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, shifting existing rows 107:205 down to 108:206
$ws.Rows("107").Insert()

# Seed the new row 107 with the same field layout as the row that is now
# directly below it (the old row 107, shifted down to 108), then overwrite
# the columns that actually carry new data for this record.
$ws.Range("A108:R108").Copy($ws.Range("A107:R107"))

$ws.Range("D107").Value = 44658
$ws.Range("J107").Value = 65
$ws.Range("K107").Value = 17000
$ws.Range("L107").Value = 17000
$ws.Range("M107").Value = 17000
$ws.Range("P107").Value = 944
